$d = $word.ActiveDocument

$d.Content.Find.Execute("2024-08-19 Monday", $true, $false, $false, $false, $false, $true, 1, $false, "2024-08-20 Tuesday", 2) | Out-Null
$d.Content.Find.Execute("97×50=", $true, $false, $false, $false, $false, $true, 1, $false, "94×11=", 2) | Out-Null
$d.Content.Find.Execute("90×16=", $true, $false, $false, $false, $false, $true, 1, $false, "26×94=", 2) | Out-Null
$d.Content.Find.Execute("36×13=", $true, $false, $false, $false, $false, $true, 1, $false, "59×63=", 2) | Out-Null
$d.Content.Find.Execute("82×64=", $true, $false, $false, $false, $false, $true, 1, $false, "22×49=", 2) | Out-Null
$d.Content.Find.Execute("88×47=", $true, $false, $false, $false, $false, $true, 1, $false, "91×71=", 2) | Out-Null
$d.Content.Find.Execute("33×49=", $true, $false, $false, $false, $false, $true, 1, $false, "42×54=", 2) | Out-Null
$d.Content.Find.Execute("46×52=", $true, $false, $false, $false, $false, $true, 1, $false, "36×31=", 2) | Out-Null
$d.Content.Find.Execute("72×43=", $true, $false, $false, $false, $false, $true, 1, $false, "78×44=", 2) | Out-Null
$d.Content.Find.Execute("49×46=", $true, $false, $false, $false, $false, $true, 1, $false, "54×62=", 2) | Out-Null
$d.Content.Find.Execute("76×74=", $true, $false, $false, $false, $false, $true, 1, $false, "46×14=", 2) | Out-Null
$d.Content.Find.Execute("25×68=", $true, $false, $false, $false, $false, $true, 1, $false, "75×67=", 2) | Out-Null
$d.Content.Find.Execute("22×44=", $true, $false, $false, $false, $false, $true, 1, $false, "32×31=", 2) | Out-Null
$d.Content.Find.Execute("44×36=", $true, $false, $false, $false, $false, $true, 1, $false, "93×93=", 2) | Out-Null
$d.Content.Find.Execute("15×22=", $true, $false, $false, $false, $false, $true, 1, $false, "69×99=", 2) | Out-Null
$d.Content.Find.Execute("58×66=", $true, $false, $false, $false, $false, $true, 1, $false, "84×71=", 2) | Out-Null
$d.Content.Find.Execute("48×97=", $true, $false, $false, $false, $false, $true, 1, $false, "49×64=", 2) | Out-Null
$d.Content.Find.Execute("98×46=", $true, $false, $false, $false, $false, $true, 1, $false, "27×55=", 2) | Out-Null
$d.Content.Find.Execute("15×13=", $true, $false, $false, $false, $false, $true, 1, $false, "13×99=", 2) | Out-Null
$d.Content.Find.Execute("13×45=", $true, $false, $false, $false, $false, $true, 1, $false, "85×45=", 2) | Out-Null
$d.Content.Find.Execute("30×27=", $true, $false, $false, $false, $false, $true, 1, $false, "54×53=", 2) | Out-Null
$d.Content.Find.Execute("11×77=", $true, $false, $false, $false, $false, $true, 1, $false, "25×75=", 2) | Out-Null
$d.Content.Find.Execute("56×38=", $true, $false, $false, $false, $false, $true, 1, $false, "98×52=", 2) | Out-Null
$d.Content.Find.Execute("85×59=", $true, $false, $false, $false, $false, $true, 1, $false, "30×82=", 2) | Out-Null
$d.Content.Find.Execute("21×75=", $true, $false, $false, $false, $false, $true, 1, $false, "25×50=", 2) | Out-Null
$d.Content.Find.Execute("99×70=", $true, $false, $false, $false, $false, $true, 1, $false, "96×19=", 2) | Out-Null
